$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2025-07-21"

$headerRange = $newSheet.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$newSheet.Range("A1").Value = "rank"
$newSheet.Range("B1").Value = "title"
$newSheet.Range("C1").Value = "author"
$newSheet.Range("D1").Value = "latest_episode"

$data = @(
    ,@(1, "宇崎ちゃんは遊びたい！", "丈(著者)", "第125話")
    ,@(2, "異世界おじさん", "殆ど死んでいる(著者)", "【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！")
    ,@(3, "時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―", "光永康則", "第６６話『六花停止』③")
    ,@(4, "生徒会にも穴はある！", "むちまろ", "第131話	ありす大ピンチ！（デジャブ編）")
    ,@(5, "転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～", "zunta(作画) はらわたさいぞう(原作)", "第30話：一秒の奪い合い③")
    ,@(6, "勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが", "絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)", "第3話 後編")
    ,@(7, "地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。", "マツモトケンゴ", "第６０話　じゃない方の戦いが始まった（２）")
    ,@(8, "帰ってください！ 阿久津さん", "長岡太一(著者)", "第192話")
    ,@(9, "元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～", "沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)", "第76話その2")
    ,@(10, "いとこのこ", "いぬちく(著者)", "休載イラスト")
    ,@(11, "異世界魔王と召喚少女の奴隷魔術", "原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大", "第126話　戦争を終わらせてみるⅡ（後編）")
    ,@(12, "実は俺、最強でした？", "原作：澄守 彩 漫画：高橋 愛", "おまけ63")
    ,@(13, "十年目、帰還を諦めた転移者はいまさら主人公になる", "原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう", "【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！")
    ,@(14, "勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～", "漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり", "第５０話　雌雄を決する器用貧乏（３）")
    ,@(15, "【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！", "島知宏 音速炒飯 有都あらゆる", "第２２食　ユクシーさんの覚悟、すごいのですわ！（３）")
    ,@(16, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第71話②")
    ,@(17, "バキ外伝 烈海王は異世界転生しても一向にかまわんッッ", "板垣恵介 猪原賽 陸井栄史", "第76話　海皇戦争")
    ,@(18, "世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜", "戸賀 環 坂木持丸 riritto", "第49話①　城のパーティーに参加してみた")
    ,@(19, "独身貴族は異世界を謳歌する ～結婚しない男の優雅なおひとりさまライフ～", "漫画：駒鳥 ひわ 原作：錬金王 キャラクター原案：三登 いつき", "第32話 独身貴族は森で写真を撮る（1）")
    ,@(20, "魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～", "漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ", "第47話 魔導具師とつながれたもの②")
    ,@(21, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第136話 よくわからないけれどスカウトされたみたいです（１）")
    ,@(22, "聖者無双", "漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime", "第90話　研究者や技術者の故郷（後半）")
    ,@(23, "アイドル辞めるけど結婚してくれますか!?", "三吉汐美(著者)", "第16話後半")
    ,@(24, "落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～", "村上よしゆき 茨木野 あるてら", "第４０話　勇者、聖女と元聖騎士と再会し、魚人を追っ払う（５）")
    ,@(25, "くらいあの子としたいこと", "碇マナツ(著者)", "第80話")
    ,@(26, "ライドンキング", "馬場康誌", "第81話 大統領と失われた神器（前編）")
    ,@(27, "小林さんちのメイドラゴン", "クール教信者", "第146話")
    ,@(28, "姫様“拷問”の時間です", "原作:春原ロビンソン　漫画:ひらけい", "拷問145")
    ,@(29, "王子様の友達", "すけろく(著者)", "【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！")
    ,@(30, "願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜", "ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)", "第4話-2：師匠と弟子の新生活")
    ,@(31, "ライブダンジョン！", "ことりりょう(作画) dy冷凍(原作) Mika Pikazo(キャラクター原案)", "第88話前半")
    ,@(32, "理想のヒモ生活", "日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)", "第86話　その1")
    ,@(33, "10年ぶりに再会したクソガキは清純美少女JKに成長していた", "緑青黒羽（漫画） 館西夕木（原作） ひげ猫（キャラクター原案）", "第5話　嫉妬、そして嫉妬（後編）")
    ,@(34, "バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～", "板垣恵介 林たかあき", "第50話 愛の試練")
    ,@(35, "俺は星間国家の悪徳領主！", "灘島かい（漫画） 三嶋与夢（原作） 高峰ナダレ（キャラクター原案）", "第39話　自惚れ")
    ,@(36, "異世界でも無難に生きたい症候群", "原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう", "第30話①")
    ,@(37, "婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版", "漫画/すたひろ 原作/Y.A", "chapter66【35話①】")
    ,@(38, "ギャルとダンジョンと周回遅れの探索英雄譚", "漫画家： 水田ケンジ 原作：榊一郎 キャラクター原案：黒獅子", "第1話")
    ,@(39, "最弱貴族に転生したので悪役たちを集めてみた", "空野進 sorani ファルまろ", "第10話　最弱貴族、部下を信じる（２）")
    ,@(40, "絶対死なないステラ姫", "光永康則 大高稲", "第１４話　絶対旅立たない（１）")
    ,@(41, "ひとりぼっちの異世界攻略", "びび（漫画） 五示正司（原作）", "第227話　自業自得です")
    ,@(42, "賢者の孫", "緒方俊輔(漫画) 吉岡剛(原作) 菊池政治(キャラクター原案)", "第94話-2")
    ,@(43, "魔法少女リリカルなのは EXCEEDS", "都築真紀 川上修一", "第４話①")
    ,@(44, "最強勇者パーティーは愛が知りたい", "山田肌襦袢", "第27話「エッチな祭りを始めたい」")
    ,@(45, "35歳独身山田、異世界村に理想のセカンドハウスを作りたい　～異世界と現実のいいとこどりライフ～", "出雲大吉(原作) 西尾洋一(作画) ゆのひと(キャラクター原案)", "第2話①")
    ,@(46, "ダメ人間の愛しかた", "岩葉(著者)", "第18話後編　ダメ人間とお姉ちゃんと彼女")
    ,@(47, "生徒会役員共", "氏家ト全", "#405")
    ,@(48, "めっちゃ召喚された件 THE COMIC", "漫画：六甲島カモメ 原作：さいとうさ キャラクター原案：ツグトク", "第46話②")
    ,@(49, "オタクに優しいギャルはいない!?", "のりしろちゃん 魚住さかな", "【#148】氷の笑み")
    ,@(50, "デスゲームに巻き込まれた山本さん、気ままにゲームバランスを崩壊させる", "ぽち(原作) カモトタツヤ(作画) 久賀フーナ(キャラクター原案)", "休載イラスト")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $item[0]
    $newSheet.Cells.Item($row, 2).Value = $item[1]
    $newSheet.Cells.Item($row, 3).Value = $item[2]
    $newSheet.Cells.Item($row, 4).Value = $item[3]
}